$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -12.0217
$ws.Range("C6").Value = -12.5004
$ws.Range("C7").Value = -13.0619
$ws.Range("E7").Value = 15.4434
$ws.Range("E12").Value = 17.7704
$ws.Range("E15").Value = 16.1936
$ws.Range("C16").Value = -14.16759999999999
$ws.Range("C20").Value = -12.2452
$ws.Range("E20").Value = 15.93779999999999
$ws.Range("E21").Value = 17.0743
$ws.Range("E22").Value = 17.0818
$ws.Range("E23").Value = 16.09779999999999
$ws.Range("C28").Value = -12.4651
$ws.Range("C29").Value = -11.6166
$ws.Range("E29").Value = 17.48030000000001
$ws.Range("C32").Value = -13.21010000000001
$ws.Range("E34").Value = 17.2696
$ws.Range("C40").Value = -12.45890000000001
$ws.Range("E42").Value = 16.36749999999999
$ws.Range("E43").Value = 17.4541
$ws.Range("E44").Value = 16.72739999999999
$ws.Range("E45").Value = 16.50729999999999
$ws.Range("C46").Value = -14.64559999999999
$ws.Range("E46").Value = 16.69620000000001
$ws.Range("E50").Value = 16.5179
$ws.Range("C51").Value = -11.4333
$ws.Range("E51").Value = 17.31150000000001
$ws.Range("C52").Value = -11.3895
$ws.Range("C57").Value = -14.23249999999998
$ws.Range("C59").Value = -13.0027
$ws.Range("C62").Value = -14.7757
$ws.Range("C66").Value = -11.2262
$ws.Range("E66").Value = 17.17980000000001
$ws.Range("E67").Value = 17.05070000000002
$ws.Range("C73").Value = -11.989
$ws.Range("C74").Value = -11.8283
$ws.Range("E79").Value = 18.12970000000002
$ws.Range("E84").Value = 16.5435
$ws.Range("C92").Value = -10.5381
$ws.Range("E92").Value = 18.06390000000002
$ws.Range("E97").Value = 16.37249999999999
$ws.Range("C100").Value = -12.5844
